# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The account-statement data table (rows 16-25) is replaced: the old set of
# workers/periods is removed and substituted with a single worker
# (PPT / 3606793 / KATIUSCA MAGALY FIGUEROA GONZALEZ) across three new
# periods (2502, 2503, 2504). The "Valor Mora" and "Cant. Trabajadores /
# Cant. Periodos" summary cells are updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old rows 18-24 (7 rows). This leaves old row 16 and 17 (the
# "regular" styled rows) in place, and shifts the old closing row 25 (which
# carries the bottom-border "last row of table" style) up to become the new
# row 18 - matching the 3-row table the target workbook ends up with.
$ws.Range("A18:A24").EntireRow.Delete() | Out-Null

# Summary fields above the table.
$ws.Range("E11").Value() = 170820
$ws.Range("C13").Value() = 1
$ws.Range("F13").Value() = 3

# New data rows (2 "regular" rows + 1 "closing" row, already carrying the
# correct formatting from the row shift above).
$ws.Range("B16").Value() = "PPT"
$ws.Range("C16").Value() = "3606793"
$ws.Range("D16").Value() = "KATIUSCA MAGALY FIGUEROA GONZALEZ"
$ws.Range("E16").Value() = "2502"
$ws.Range("F16").Value() = 56940

$ws.Range("B17").Value() = "PPT"
$ws.Range("C17").Value() = "3606793"
$ws.Range("D17").Value() = "KATIUSCA MAGALY FIGUEROA GONZALEZ"
$ws.Range("E17").Value() = "2503"
$ws.Range("F17").Value() = 56940

$ws.Range("B18").Value() = "PPT"
$ws.Range("C18").Value() = "3606793"
$ws.Range("D18").Value() = "KATIUSCA MAGALY FIGUEROA GONZALEZ"
$ws.Range("E18").Value() = "2504"
$ws.Range("F18").Value() = 56940
